# Daily attendance processing - normalize "Modified By" (column G) ordering
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Known reorderings of the comma-separated "modified by" list observed in the
# attendance export. Only cells whose text matches exactly one of these
# values (left-hand side) get rewritten to the normalized ordering
# (right-hand side); everything else is left untouched.
$map = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
    "backup@backdoor.com, System, system" = "backup@backdoor.com, system, System"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G$r")
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
